$wb = $excel.ActiveWorkbook

# Rename the first worksheet from "ERP_vintages" to "ERP" for consistency
# with the workbook/file naming.
$ws = $wb.Worksheets.Item(1)
$ws.Name = "ERP"
